# Normalize the "Recorded By" column (column G) on the "Session Analysis
# Results" sheet so that the "System" entry (added automatically by the
# attendance system) is listed first, ahead of the human/user entries.
#
# Examples of the transform applied:
#   "dnasr281@gmail.com, System"            -> "System, dnasr281@gmail.com"
#   "admin@admin.com, System"               -> "System, admin@admin.com"
#   "backup@backdoor.com, system, System"   -> "backup@backdoor.com, System, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Value2

    if ($value -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($value -eq "admin@admin.com, System") {
        $cell.Value = "System, admin@admin.com"
    }
    elseif ($value -eq "backup@backdoor.com, system, System") {
        $cell.Value = "backup@backdoor.com, System, system"
    }
}
